$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cr = [char]13
$rsquo = [char]0x2019

# --- Row 2 (David Maxson), Column 4 (R3/R4): was a single empty
#     paragraph; becomes five paragraphs of text. Joining the text with
#     carriage returns in one assignment lets Word create the extra
#     paragraphs while inheriting the original (pre-existing) run
#     formatting for every new paragraph.
$cellR3R4 = $t.Rows.Item(2).Cells.Item(4)
$textR3R4 = "sys_call_irq" + $cr + "User Manual" + $cr + "Programmer" + $rsquo + "s Manual" + $cr + "Idle" + $cr + "infinite cmd"
$cellR3R4.Range.Text = $textR3R4

# --- Row 2 (David Maxson), Column 5 (R5): was a single empty paragraph;
#     becomes two empty paragraphs (still no text). Duplicate the
#     paragraph via InsertParagraphAfter so no stray <w:t/> is created.
$cellR5 = $t.Rows.Item(2).Cells.Item(5)
$cellR5.Range.Paragraphs.Item(1).Range.InsertParagraphAfter()

# --- Font change: every cell in the R3/R4, R5, R6 columns (columns 4-6)
#     for all four member rows (rows 2-5) switches from "Times New Roman"
#     to "arial". This covers the two cells edited above plus the nine
#     other cells that only change font (no content change).
for ($r = 2; $r -le 5; $r++) {
    for ($c = 4; $c -le 6; $c++) {
        $cell = $t.Rows.Item($r).Cells.Item($c)
        $cell.Range.Font.Name = "arial"
    }
}
